$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two re-dated report labels in the header row (sharedStrings text)
$ws.Range("I9").Value = "1402-03-07 (8)"
$ws.Range("M9").Value = "1402-03-07 (2)"

$cols = @("D","E","F","G","H","I","J","K","L","M")

$row11 = @(1788320, 1648559, 1644247, 1804624, 2081045, 1702161, 1871508, 3044623, 3786661, 3652655)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $row11[$i]
}

$row12 = @(-674879, -792007, -778096, -1120104, -1097998, -1187928, -1121670, -1224797, -1918305, -2014692)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $row12[$i]
}

$row13 = @(1113441, 856552, 866151, 684520, 983047, 514233, 749838, 1819826, 1868356, 1637963)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $row13[$i]
}

$row14 = @(-41136, -68244, -63510, -84635, -84818, -204524, -125299, -206716, -47358, -275023)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $row14[$i]
}

$row15 = @(0, 0, 0, 0, 0, 0, 0, -23444, 0, 23444)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Value = $row15[$i]
}

$row16 = @(-12377, 2548, -17405, 11135, -16371, -3115, 7642, 23985, -72658, 116366)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $row16[$i]
}

$row17 = @(1059928, 790856, 785236, 611020, 881858, 306594, 632181, 1613651, 1748340, 1502750)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
}

$row18 = @(-76682, -86855, -98675, -133136, -160060, -166741, -196809, -267150, -287934, -259743)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "18").Value = $row18[$i]
}

$row19 = @(697, 2010, 127348, 10528, 1566, 1663, 304423, -99943, -102685, -10265)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
}

$row20 = @(983943, 706011, 813909, 488412, 723364, 141516, 739795, 1246558, 1357721, 1232742)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
}

$row21 = @(-261378, 23674, -183130, -76847, -162490, 59709, -97961, -272057, -382398, 1626)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "21").Value = $row21[$i]
}

$row22 = @(722565, 729685, 630779, 411565, 560874, 201225, 641834, 974501, 975323, 1234368)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $row22[$i]
}

$row23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "23").Value = $row23[$i]
}

$row24 = @(722565, 729685, 630779, 411565, 560874, 201225, 641834, 974501, 975323, 1234368)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "24").Value = $row24[$i]
}

$row25 = @(702, 477, 474, 269, 367, 83, 264, 401, 269, 340)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "25").Value = $row25[$i]
}

$row26 = @(1030000, 1530000, 1330000, 1530000, 1530000, 2430000, 2430000, 2430000, 3630000, 3630000)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
}

$row27 = @(199, 201, 174, 113, 155, 55, 177, 268, 269, 340)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
}
